# Update the descriptive-statistics table with the revised (publication) figures.
# The new table grows from rows 2-23 to rows 2-32: rows 2-16 keep their
# measure/statistic labels but get new values; rows 17-21 are relabeled from
# "infpernode" to "inf_peak" with new values; rows 22-26 introduce a new
# "time_peak" block; rows 27-31 re-introduce the "infpernode" block (with the
# old mean/max/min figures shifted accordingly and new median/sd); and the
# final "never_infected" / total_pct row moves from row 22 down to row 32
# with an updated value. "max_infected_nonepi" (old row 23) is removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{R=2;  A="duration";            B="min";       C=33},
    @{R=3;  A="duration";            B="mean";      C=37.2},
    @{R=4;  A="duration";            B="median";    C=37},
    @{R=5;  A="duration";            B="max";       C=43},
    @{R=6;  A="duration";            B="sd";        C=2.12251042519242},
    @{R=7;  A="infections";          B="min";       C=68.86011820637988},
    @{R=8;  A="infections";          B="mean";      C=69.45028876784745},
    @{R=9;  A="infections";          B="median";    C=69.4609857831491},
    @{R=10; A="infections";          B="max";       C=70.09794934085046},
    @{R=11; A="infections";          B="sd";        C=0.2586293839891742},
    @{R=12; A="tti_mean";            B="min";       C=12.25874554004724},
    @{R=13; A="tti_mean";            B="mean";      C=13.92888260983688},
    @{R=14; A="tti_mean";            B="median";    C=13.84910333471393},
    @{R=15; A="tti_mean";            B="max";       C=16.31782648724041},
    @{R=16; A="tti_mean";            B="sd";        C=0.7374389936342706},
    @{R=17; A="inf_peak";            B="min";       C=12.23979852330211},
    @{R=18; A="inf_peak";            B="mean";      C=12.96476958903796},
    @{R=19; A="inf_peak";            B="median";    C=12.96675201584923},
    @{R=20; A="inf_peak";            B="max";       C=13.87202206261515},
    @{R=21; A="inf_peak";            B="sd";        C=0.3376795507266901},
    @{R=22; A="time_peak";           B="min";       C=13},
    @{R=23; A="time_peak";           B="mean";      C=14.48},
    @{R=24; A="time_peak";           B="median";    C=14},
    @{R=25; A="time_peak";           B="max";       C=17},
    @{R=26; A="time_peak";           B="sd";        C=0.8466213302106655},
    @{R=27; A="infpernode";          B="min";       C=0},
    @{R=28; A="infpernode";          B="mean";      C=69.45028876784745},
    @{R=29; A="infpernode";          B="median";    C=82},
    @{R=30; A="infpernode";          B="max";       C=100},
    @{R=31; A="infpernode";          B="sd";        C=33.15236740884213},
    @{R=32; A="never_infected";      B="total_pct"; C=0.2790276621263397}
)

foreach ($row in $rows) {
    $ws.Cells.Item($row.R, 1).Value = $row.A
    $ws.Cells.Item($row.R, 2).Value = $row.B
    $ws.Cells.Item($row.R, 3).Value = $row.C
}
